$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.394.85'
$ws.Range('E2').Value = '  -2.34%  '
$ws.Range('D3').Value = '3.692.66'
$ws.Range('E3').Value = '  -3.07%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = "'691.64"
$ws.Range('E5').Value = '  -1.38%  '
$ws.Range('D6').Value = "'162.40"
$ws.Range('E6').Value = '  -5.62%  '
$ws.Range('D7').Value = '3.691.35'
$ws.Range('E7').Value = '  -3.12%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').Value = "'0.498"
$ws.Range('E9').Value = '  -4.88%  '
$ws.Range('E10').Value = '  -8.26%  '
$ws.Range('D11').Value = "'7.40"
$ws.Range('E11').Value = '  -1.65%  '
$ws.Range('D12').Value = "'0.441"
$ws.Range('E12').Value = '  -5.17%  '
$ws.Range('D13').Value = "'0.0000239"
$ws.Range('E13').Value = '  -5.54%  '
$ws.Range('D14').Value = "'33.35"
$ws.Range('E14').Value = '  -7.33%  '
$ws.Range('D15').Value = '4.315.56'
$ws.Range('E15').Value = '  -3.13%  '
$ws.Range('D16').Value = '3.692.16'
$ws.Range('E16').Value = '  -3.50%  '
$ws.Range('D17').Value = '69.434.02'
$ws.Range('E17').Value = '  -2.37%  '
$ws.Range('E18').Value = '  -0.61%  '
$ws.Range('D19').Value = "'16.18"
$ws.Range('E19').Value = '  -7.45%  '
$ws.Range('D20').Value = "'6.57"
$ws.Range('E20').Value = '  -7.92%  '
$ws.Range('D21').Value = "'480.61"
$ws.Range('E21').Value = '  -6.21%  '
$ws.Range('D22').Value = "'9.98"
$ws.Range('D23').Value = "'0.662"
$ws.Range('E23').Value = '  -7.32%  '
$ws.Range('D24').Value = "'79.85"
$ws.Range('E24').Value = '  -4.84%  '
$ws.Range('D25').Value = '3.839.65'
$ws.Range('E25').Value = '  -3.10%  '
$ws.Range('E26').Value = '  -9.69%  '
$ws.Range('E27').Value = '  +0.05%  '
$ws.Range('D28').Value = "'11.32"
$ws.Range('E28').Value = '  -6.16%  '
$ws.Range('D29').Value = "'9.48"
$ws.Range('E29').Value = '  -8.88%  '
$ws.Range('E30').Value = '  -10.64%  '
$ws.Range('E31').Value = '  -10.03%  '
$ws.Range('E32').Value = '  -7.99%  '
$ws.Range('B34').Value = 'Kaspa'
$ws.Range('C34').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D34').Value = "'0.167"
$ws.Range('E34').Value = '  -5.50%  '
$ws.Range('B35').Value = 'Binance-PegBSC-USD'
$ws.Range('C35').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D35').Value = "'1.00"
$ws.Range('E35').Value = '  +0.24%  '
$ws.Range('B36').Value = 'EthereumClassic'
$ws.Range('C36').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D36').Value = "'26.93"
$ws.Range('E36').Value = '  -7.23%  '
$ws.Range('D37').Value = '3.666.26'
$ws.Range('E37').Value = '  -2.81%  '
$ws.Range('D38').Value = "'8.45"
$ws.Range('E38').Value = '  -7.59%  '
$ws.Range('D39').Value = "'6.30"
$ws.Range('E39').Value = '  +4.99%  '
$ws.Range('D40').Value = "'2.32"
$ws.Range('E40').Value = '  -2.45%  '
$ws.Range('D41').Value = "'0.0928"
$ws.Range('E41').Value = '  -8.05%  '
$ws.Range('E43').Value = '  -0.03%  '
$ws.Range('E44').Value = '  -6.37%  '
$ws.Range('D45').Value = "'164.52"
$ws.Range('E45').Value = '  -5.02%  '
$ws.Range('D46').Value = "'48.02"
$ws.Range('E46').Value = '  -2.69%  '
$ws.Range('D47').Value = "'30.09"
$ws.Range('E47').Value = '  +2.94%  '
$ws.Range('D48').Value = "'2.79"
$ws.Range('E48').Value = '  -15.74%  '
$ws.Range('B49').Value = 'SuiNetwork'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D49').Value = "'1.14"
$ws.Range('E49').Value = '  -1.09%  '
$ws.Range('B50').Value = 'ONDO'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D50').Value = "'1.34"
$ws.Range('E50').Value = '  -2.35%  '
$ws.Range('D51').Value = "'0.000283"
$ws.Range('E51').Value = '  -10.08%  '
